# Apply the "Dep Ed Closures" daily-update edit:
#  - Update the "On this page" summary text (row 7) with the new date/time stamp
#  - Remove the resolved school closure (St John's School, FRANKSTON EAST) and the
#    old merged TAFE/region row (rows 79 & 80), which shifts every subsequent row
#    up by two and collapses the "Schools closed" / "TAFE" blurb for the
#    South-Eastern Victoria region into the new "no closures" wording (row 78)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "On this page" summary (date/time stamp refreshed)
$ws.Range("A7").Value = "On this pageCurrent school and early childhood service, TAFE closures and relocations:Bus service cancellations or alterationsCurrent school and early childhood service, TAFE closures and relocations for Friday 2 October, (as at 9:50am, 2 October)South-Eastern Victoria RegionEarly childhood services"

# 2) Delete the two rows that no longer apply:
#    row 79 - "li: St John's School, FRANKSTON EAST" (closure resolved)
#    row 80 - "TAFE...North-Eastern Victoria Region..." merged row
# Deleting row 79 twice removes both rows and shifts everything below up by two.
$ws.Rows.Item(79).Delete()
$ws.Rows.Item(79).Delete()

# 3) Update the now-merged "Schools closed" / "TAFE" summary for the
#    South-Eastern Victoria region (row 78) to reflect no closures reported
$ws.Range("A78").Value = "Schools closedThe Department hasnotbeen advised of any school closures.TAFETheDepartment hasnotbeen advised of any TAFE closures.North-Eastern Victoria RegionEarly childhood services"
